# Adjust Investment Summary table column widths for better formatting
# (also touches the Timeline & Milestones table on slide 5, which shares
# the same table style / layout placeholder sizing).

$p = $ppt.ActivePresentation

# --- Slide 5: "Timeline & Milestones" table (4 columns) ---
$s5 = $p.Slides.Item(5)
$tbl5 = $s5.Shapes.Item(3).Table
$tbl5.Columns.Item(1).Width = 871093 / 12700
$tbl5.Columns.Item(3).Width = 1306639 / 12700
$tbl5.Columns.Item(4).Width = 4355466 / 12700

# --- Slide 8: "Investment Summary" table (7 columns) ---
$s8 = $p.Slides.Item(8)
$tbl8 = $s8.Shapes.Item(3).Table
$tbl8.Columns.Item(1).Width = 1742186 / 12700
$tbl8.Columns.Item(2).Width = 1045311 / 12700
$tbl8.Columns.Item(3).Width = 2003514 / 12700
$tbl8.Columns.Item(4).Width = 1132421 / 12700
$tbl8.Columns.Item(5).Width = 871093 / 12700
$tbl8.Columns.Item(6).Width = 871093 / 12700
$tbl8.Columns.Item(7).Width = 1045311 / 12700
